$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 349, pushing existing rows 349-442 down to 351-444
$ws.Rows.Item(349).Resize(2).Insert()

# Row 349: new "Primera" record dated 44964
$ws.Cells.Item(349, 1).Value = 11
$ws.Cells.Item(349, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(349, 3).Value = "Bíobío"
$ws.Cells.Item(349, 4).Value = 44964
$ws.Cells.Item(349, 5).Value = 8
$ws.Cells.Item(349, 6).Value = 100114014
$ws.Cells.Item(349, 7).Value = "Betarraga"
$ws.Cells.Item(349, 8).Value = "Sin especificar"
$ws.Cells.Item(349, 9).Value = "Primera"
$ws.Cells.Item(349, 10).Value = 500
$ws.Cells.Item(349, 11).Value = 700
$ws.Cells.Item(349, 12).Value = 800
$ws.Cells.Item(349, 13).Value = 740
$ws.Cells.Item(349, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(349, 15).Value = "Región Metropolitana"
$ws.Cells.Item(349, 16).Value = 148
$ws.Cells.Item(349, 17).Value = 5
$ws.Cells.Item(349, 18).Value = "Hortaliza"

# Row 350: new "Segunda" record dated 44964
$ws.Cells.Item(350, 1).Value = 11
$ws.Cells.Item(350, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(350, 3).Value = "Bíobío"
$ws.Cells.Item(350, 4).Value = 44964
$ws.Cells.Item(350, 5).Value = 8
$ws.Cells.Item(350, 6).Value = 100114014
$ws.Cells.Item(350, 7).Value = "Betarraga"
$ws.Cells.Item(350, 8).Value = "Sin especificar"
$ws.Cells.Item(350, 9).Value = "Segunda"
$ws.Cells.Item(350, 10).Value = 300
$ws.Cells.Item(350, 11).Value = 600
$ws.Cells.Item(350, 12).Value = 600
$ws.Cells.Item(350, 13).Value = 600
$ws.Cells.Item(350, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(350, 15).Value = "Región Metropolitana"
$ws.Cells.Item(350, 16).Value = 120
$ws.Cells.Item(350, 17).Value = 5
$ws.Cells.Item(350, 18).Value = "Hortaliza"
